$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'26.110.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '

$ws.Range("D3").Value = "`'1.647.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.71%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = "`'218.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '

$ws.Range("D6").Value = "`'0.5203"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.71%  '

$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").Value = "`'0.2619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.03%  '

$ws.Range("D9").Value = "`'0.06299"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").Value = "`'20.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.82%  '

$ws.Range("D11").Value = "`'0.07642"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.85%  '

$ws.Range("D12").Value = "`'4.583"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.53%  '

$ws.Range("D13").Value = "`'1.662.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.48%  '

$ws.Range("D14").Value = "`'1.874.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.81%  '

$ws.Range("D15").Value = "`'0.5569"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.03%  '

$ws.Range("D16").Value = "`'0.0₅8115"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.60%  '

$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("D18").Value = "`'26.061.55"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("D20").Value = "`'4.587"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.93%  '

$ws.Range("D21").Value = "`'194.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").Value = "`'10.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.95%  '

$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").Value = "`'145.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.22%  '

$ws.Range("D26").Value = "`'0.1181"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.51%  '

$ws.Range("D27").Value = "`'7.192"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.50%  '

$ws.Range("D28").Value = "`'1.536"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.39%  '

$ws.Range("D29").Value = "`'15.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.38%  '

$ws.Range("D30").Value = "`'0.05441"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.59%  '

$ws.Range("E31").Value = '  -0.70%  '

$ws.Range("D32").Value = "`'3.432"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.73%  '

$ws.Range("D33").Value = "`'3.324"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.45%  '

$ws.Range("D34").Value = "`'1.558"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.79%  '

$ws.Range("D35").Value = "`'2.412"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.34%  '

$ws.Range("D36").Value = "`'2.781"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.07%  '

$ws.Range("D37").Value = "`'0.9425"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.26%  '

$ws.Range("D38").Value = "`'0.5591"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.76%  '

$ws.Range("E39").Value = '  -2.34%  '

$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("D41").Value = "`'5.733"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.23%  '

$ws.Range("D42").Value = "`'1.027.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.18%  '

$ws.Range("D43").Value = "`'0.8198"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.28%  '

$ws.Range("D44").Value = "`'100.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.40%  '

$ws.Range("D45").Value = "`'1.785.53"
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = '  +7.50%  '

$ws.Range("D47").Value = "`'57.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.04%  '

$ws.Range("D48").Value = "`'0.9993"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("D49").Value = "`'0.4317"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.86%  '

$ws.Range("D50").Value = "`'7.871"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.90%  '

$ws.Range("D51").Value = "`'0.05098"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.21%  '

